$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Update rows 10-17: column A -> 1, column C -> 6
for ($r = 10; $r -le 17; $r++) {
    $ws1.Cells.Item($r, 1).Value = 1
    $ws1.Cells.Item($r, 3).Value = 6
}

# Row 10 (first block) F:I mark-up values reset to 0
$ws1.Range("F10:I10").Value = 0

# Delete rows 18-25 entirely (shifts cells up)
[void]$ws1.Range("A18:A25").EntireRow.Delete()

[void]$wb.Worksheets.Item("Sheet2").Delete()

[void]$ws1.Range("D26").Select()

